$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Yes"
$ws.Range("F2").Value = "Na"
$ws.Range("J2").Value = 1

$ws.Range("B3").Value = "Yes"
$ws.Range("F3").Value = "Na"
$ws.Range("J3").Value = 1
